$wb = $excel.ActiveWorkbook
$tc = $wb.Worksheets.Item("TestCases")

# Add ObjectRepository sheet after TestCases
$repo = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $tc)
$repo.Name = "ObjectRepository"

# Populate ObjectRepository sheet
$repo.Range("A1").Value = "Alias"
$repo.Range("B1").Value = "Identifier"
$repo.Range("C1").Value = "Locator"

$repo.Range("A2").Value = "SigninButton"
$repo.Range("B2").Value = "xpath"
$repo.Range("C2").Value = "//button[@id='login']"

$repo.Range("A3").Value = "EmailInput"
$repo.Range("B3").Value = "xpath"
$repo.Range("C3").Value = "//input[@id='email']"

$repo.Range("A4").Value = "Next"
$repo.Range("B4").Value = "xpath"
$repo.Range("C4").Value = "//button[@id='next']"

$repo.Range("A5").Value = "EmailPassword"
$repo.Range("B5").Value = "xpath"
$repo.Range("C5").Value = "//input[@id='password']"

$repo.Range("A6").Value = "LoginButton"
$repo.Range("B6").Value = "xpath"
$repo.Range("C6").Value = "//button[@id='taLogin']"

# Match styling used across the rest of the workbook (same font style as TestCases sheet)
$tc.Range("A1").Copy()
$repo.Range("A1:C6").PasteSpecial(-4122)  # xlPasteFormats

# Update TestCases sheet rows 13-17: remove column E, replace F values with alias names
$tc.Range("E13").Clear()
$tc.Range("F13").Value = "SigninButton"

$tc.Range("E14").Clear()
$tc.Range("F14").Value = "EmailInput"

$tc.Range("E15").Clear()
$tc.Range("F15").Value = "Next"

$tc.Range("E16").Clear()
$tc.Range("F16").Value = "EmailPassword"

$tc.Range("E17").Clear()
$tc.Range("F17").Value = "LoginButton"

Write-Output "done"
